# Auto-generated edit script: updates H:N profit-analysis columns
# across all 8 sheets per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 37.545456
$ws.Range("I5").Value = 43.125
$ws.Range("K5").Value = 43.125
$ws.Range("M5").Value = 71.875
$ws.Range("H9").Value = 108.5
$ws.Range("I9").Value = 294.33334
$ws.Range("J9").Value = 28.857143
$ws.Range("K9").Value = 294.33334
$ws.Range("L9").Value = 28.857143
$ws.Range("M9").Value = -125.33334
$ws.Range("N9").Value = -366.857143
$ws.Range("H31").Value = 831
$ws.Range("I31").Value = 831
$ws.Range("K31").Value = 2493
$ws.Range("M31").Value = -2263
$ws.Range("H33").Value = 30481
$ws.Range("I33").Value = 32135.938
$ws.Range("J33").Value = 4002
$ws.Range("K33").Value = 32135.938
$ws.Range("L33").Value = 4002
$ws.Range("M33").Value = -31906.938
$ws.Range("N33").Value = -4460
$ws.Range("H40").Value = 9399.429
$ws.Range("I40").Value = 8930.666999999999
$ws.Range("J40").Value = 9751
$ws.Range("K40").Value = 8930.666999999999
$ws.Range("L40").Value = 9751
$ws.Range("M40").Value = -8755.666999999999
$ws.Range("N40").Value = -10101
$ws.Range("H41").Value = 1193
$ws.Range("I41").Value = 1969.4
$ws.Range("J41").Value = 761.6667
$ws.Range("K41").Value = 1969.4
$ws.Range("L41").Value = 761.6667
$ws.Range("M41").Value = -1529.4
$ws.Range("N41").Value = -1641.6667
$ws.Range("H42").Value = 1806.8
$ws.Range("I42").Value = 188.6
$ws.Range("J42").Value = 3425
$ws.Range("K42").Value = 565.8
$ws.Range("L42").Value = 10275
$ws.Range("M42").Value = -335.8
$ws.Range("N42").Value = -10735
$ws.Range("H43").Value = 5049.9463
$ws.Range("J43").Value = 4734.533
$ws.Range("L43").Value = 4734.533
$ws.Range("N43").Value = -4872.533
$ws.Range("H51").Value = 13772.5
$ws.Range("J51").Value = 3597.1428
$ws.Range("L51").Value = 3597.1428
$ws.Range("N51").Value = -4565.1428
$ws.Range("H53").Value = 1078.96
$ws.Range("J53").Value = 823.0714
$ws.Range("L53").Value = 823.0714
$ws.Range("N53").Value = -2097.0714
$ws.Range("H61").Value = 10000
$ws.Range("J61").Value = 10000
$ws.Range("L61").Value = 30000
$ws.Range("N61").Value = -30344
$ws.Range("H62").Value = 4999.5
$ws.Range("I62").Value = 4999
$ws.Range("K62").Value = 4999
$ws.Range("M62").Value = -4375
$ws.Range("H65").Value = 4999.5
$ws.Range("I65").Value = 4999
$ws.Range("K65").Value = 24995
$ws.Range("M65").Value = -21875
$ws.Range("H70").Value = 3436.125
$ws.Range("I70").Value = 1499.5
$ws.Range("J70").Value = 4081.6667
$ws.Range("K70").Value = 4498.5
$ws.Range("L70").Value = 12245.0001
$ws.Range("M70").Value = -4228.5
$ws.Range("N70").Value = -12785.0001
$ws.Range("H73").Value = 3436.125
$ws.Range("I73").Value = 1499.5
$ws.Range("J73").Value = 4081.6667
$ws.Range("K73").Value = 4498.5
$ws.Range("L73").Value = 12245.0001
$ws.Range("M73").Value = -3562.5
$ws.Range("N73").Value = -14117.0001
$ws.Range("H76").Value = 3915.4167
$ws.Range("J76").Value = 4999.5
$ws.Range("L76").Value = 4999.5
$ws.Range("N76").Value = -5629.5
$ws.Range("H79").Value = 3915.4167
$ws.Range("J79").Value = 4999.5
$ws.Range("L79").Value = 4999.5
$ws.Range("N79").Value = -7183.5
$ws.Range("H86").Value = 3003.3684
$ws.Range("I86").Value = 2173.6
$ws.Range("J86").Value = 3925.3333
$ws.Range("K86").Value = 2173.6
$ws.Range("L86").Value = 3925.3333
$ws.Range("M86").Value = -1050.6
$ws.Range("N86").Value = -6171.3333
$ws.Range("H89").Value = 3003.3684
$ws.Range("I89").Value = 2173.6
$ws.Range("J89").Value = 3925.3333
$ws.Range("K89").Value = 10868
$ws.Range("L89").Value = 19626.6665
$ws.Range("M89").Value = -5252
$ws.Range("N89").Value = -30858.6665
$ws.Range("H98").Value = 1623.6154
$ws.Range("J98").Value = 1824.5
$ws.Range("L98").Value = 1824.5
$ws.Range("N98").Value = -4820.5
$ws.Range("H103").Value = 4933.4287
$ws.Range("I103").Value = 10614.667
$ws.Range("K103").Value = 31844.001
$ws.Range("M103").Value = -31258.001
$ws.Range("H111").Value = 1286.25
$ws.Range("I111").Value = 1286.25
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3858.75
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -791.75
$ws.Range("H115").Value = 295.625
$ws.Range("I115").Value = 295.625
$ws.Range("K115").Value = 886.875
$ws.Range("M115").Value = 680.125
$ws.Range("H118").Value = 1247.8182
$ws.Range("I118").Value = 387
$ws.Range("J118").Value = 5121.5
$ws.Range("K118").Value = 1161
$ws.Range("L118").Value = 15364.5
$ws.Range("M118").Value = 496
$ws.Range("N118").Value = -18678.5
$ws.Range("H122").Value = 1623.6154
$ws.Range("J122").Value = 1824.5
$ws.Range("L122").Value = 5473.5
$ws.Range("N122").Value = -10373.5
$ws.Range("H125").Value = 6874
$ws.Range("I125").Value = 6248.8
$ws.Range("K125").Value = 56239.2
$ws.Range("M125").Value = -53779.2
$ws.Range("H127").Value = 1607.875
$ws.Range("I127").Value = 963
$ws.Range("K127").Value = 2889
$ws.Range("M127").Value = 2071
$ws.Range("H129").Value = 2405.7896
$ws.Range("I129").Value = 482
$ws.Range("K129").Value = 1446
$ws.Range("M129").Value = 3554
$ws.Range("H132").Value = 47560.047
$ws.Range("I132").Value = 49763.43
$ws.Range("J132").Value = 1289
$ws.Range("K132").Value = 149290.29
$ws.Range("L132").Value = 3867
$ws.Range("M132").Value = -146760.29
$ws.Range("N132").Value = -8927
$ws.Range("H135").Value = 1182.9
$ws.Range("I135").Value = 1238.7778
$ws.Range("K135").Value = 11149.0002
$ws.Range("M135").Value = -8614.0002
$ws.Range("H137").Value = 450000000
$ws.Range("I137").Value = 450000000
$ws.Range("K137").Value = 1350000000
$ws.Range("M137").Value = -1349997450
$ws.Range("H138").Value = 3167.111
$ws.Range("J138").Value = 2800.535
$ws.Range("L138").Value = 8401.605
$ws.Range("N138").Value = -18681.605
$ws.Range("H141").Value = 4094.625
$ws.Range("I141").Value = 4094.625
$ws.Range("K141").Value = 12283.875
$ws.Range("M141").Value = -7103.875
$ws.Range("N111").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3240.3
$ws.Range("I32").Value = 3330.3572
$ws.Range("K32").Value = 3330.3572
$ws.Range("M32").Value = -3043.3572
$ws.Range("H45").Value = 1382.2
$ws.Range("I45").Value = 1421.0834
$ws.Range("K45").Value = 1421.0834
$ws.Range("M45").Value = -1044.0834
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16450
$ws.Range("H61").Value = 2499
$ws.Range("I61").Value = 2499
$ws.Range("K61").Value = 2499
$ws.Range("M61").Value = -2287
$ws.Range("H74").Value = 5559657
$ws.Range("I74").Value = 5127.5
$ws.Range("J74").Value = 27777776
$ws.Range("K74").Value = 5127.5
$ws.Range("L74").Value = 27777776
$ws.Range("M74").Value = -4253.5
$ws.Range("N74").Value = -27779524
$ws.Range("H77").Value = 5559657
$ws.Range("I77").Value = 5127.5
$ws.Range("J77").Value = 27777776
$ws.Range("K77").Value = 25637.5
$ws.Range("L77").Value = 138888880
$ws.Range("M77").Value = -21269.5
$ws.Range("N77").Value = -138897616
$ws.Range("H93").Value = 45052.168
$ws.Range("I93").Value = 35400
$ws.Range("J93").Value = 49878.25
$ws.Range("K93").Value = 35400
$ws.Range("L93").Value = 49878.25
$ws.Range("M93").Value = -32904
$ws.Range("N93").Value = -54870.25
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H122").Value = 8728.5
$ws.Range("I122").Value = 8000
$ws.Range("K122").Value = 24000
$ws.Range("M122").Value = -21550
$ws.Range("I132").Value = 2496
$ws.Range("K132").Value = 7488
$ws.Range("M132").Value = -4958
$ws.Range("H136").Value = 2499
$ws.Range("I136").Value = 2499
$ws.Range("K136").Value = 7497
$ws.Range("M136").Value = -4947
$ws.Range("N106").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 366.42856
$ws.Range("I22").Value = 356.30768
$ws.Range("J22").Value = 382.875
$ws.Range("K22").Value = 356.30768
$ws.Range("L22").Value = 382.875
$ws.Range("M22").Value = -183.30768
$ws.Range("N22").Value = -728.875
$ws.Range("H94").Value = 2935.5715
$ws.Range("I94").Value = 2735.5715
$ws.Range("J94").Value = 3135.5715
$ws.Range("K94").Value = 2735.5715
$ws.Range("L94").Value = 3135.5715
$ws.Range("M94").Value = -2284.5715
$ws.Range("N94").Value = -4037.5715
$ws.Range("H99").Value = 1917.4
$ws.Range("I99").Value = 1960
$ws.Range("J99").Value = 1889
$ws.Range("K99").Value = 1960
$ws.Range("L99").Value = 1889
$ws.Range("M99").Value = -462
$ws.Range("N99").Value = -4885
$ws.Range("H134").Value = 13891019
$ws.Range("I134").Value = 17859062
$ws.Range("J134").Value = 2868
$ws.Range("K134").Value = 53577186
$ws.Range("L134").Value = 8604
$ws.Range("M134").Value = -53574651
$ws.Range("N134").Value = -13674

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1404.7
$ws.Range("I16").Value = 1394.1111
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1394.1111
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1107.1111
$ws.Range("N16").Value = -2074
$ws.Range("H31").Value = 2745
$ws.Range("I31").Value = 2745
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2745
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2450
$ws.Range("H34").Value = 2745
$ws.Range("I34").Value = 2745
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2745
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2543
$ws.Range("H62").Value = 14533.833
$ws.Range("I62").Value = 3300
$ws.Range("K62").Value = 3300
$ws.Range("M62").Value = -2676
$ws.Range("H65").Value = 14533.833
$ws.Range("I65").Value = 3300
$ws.Range("K65").Value = 16500
$ws.Range("M65").Value = -13380
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H86").Value = 5980.6665
$ws.Range("I86").Value = 5632.6665
$ws.Range("J86").Value = 6328.6665
$ws.Range("K86").Value = 5632.6665
$ws.Range("L86").Value = 6328.6665
$ws.Range("M86").Value = -4509.6665
$ws.Range("N86").Value = -8574.666499999999
$ws.Range("H89").Value = 5980.6665
$ws.Range("I89").Value = 5632.6665
$ws.Range("J89").Value = 6328.6665
$ws.Range("K89").Value = 28163.3325
$ws.Range("L89").Value = 31643.3325
$ws.Range("M89").Value = -22547.3325
$ws.Range("N89").Value = -42875.3325
$ws.Range("H93").Value = 10515.909
$ws.Range("J93").Value = 14193.667
$ws.Range("L93").Value = 14193.667
$ws.Range("N93").Value = -17937.667
$ws.Range("H105").Value = 3504.4546
$ws.Range("I105").Value = 3147.375
$ws.Range("J105").Value = 3708.5
$ws.Range("K105").Value = 3147.375
$ws.Range("L105").Value = 3708.5
$ws.Range("M105").Value = -1400.375
$ws.Range("N105").Value = -7202.5
$ws.Range("H113").Value = 1404.7
$ws.Range("I113").Value = 1394.1111
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1394.1111
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 775.8888999999999
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 14842.947
$ws.Range("I122").Value = 1626.5385
$ws.Range("J122").Value = 43478.5
$ws.Range("K122").Value = 4879.6155
$ws.Range("L122").Value = 130435.5
$ws.Range("M122").Value = -2429.6155
$ws.Range("N122").Value = -135335.5
$ws.Range("H132").Value = 4814.375
$ws.Range("I132").Value = 5405.2
$ws.Range("K132").Value = 16215.6
$ws.Range("M132").Value = -13685.6
$ws.Range("H134").Value = 2786.2856
$ws.Range("I134").Value = 2786.2856
$ws.Range("K134").Value = 8358.856800000001
$ws.Range("M134").Value = -5823.856800000001
$ws.Range("N31").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 933.44446
$ws.Range("I5").Value = 933.44446
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2800.33338
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2688.33338
$ws.Range("H8").Value = 728.36365
$ws.Range("I8").Value = 728.36365
$ws.Range("K8").Value = 2185.09095
$ws.Range("M8").Value = -2046.09095
$ws.Range("H32").Value = 2428.5715
$ws.Range("J32").Value = 2428.5715
$ws.Range("L32").Value = 7285.7145
$ws.Range("N32").Value = -7851.7145
$ws.Range("H34").Value = 4282.4644
$ws.Range("I34").Value = 786.7778
$ws.Range("K34").Value = 2360.3334
$ws.Range("M34").Value = -2276.3334
$ws.Range("H46").Value = 25003584
$ws.Range("I46").Value = 50000176
$ws.Range("K46").Value = 150000528
$ws.Range("M46").Value = -150000437
$ws.Range("H98").Value = 439.22223
$ws.Range("I98").Value = 359.66666
$ws.Range("J98").Value = 479
$ws.Range("K98").Value = 1078.99998
$ws.Range("L98").Value = 1437
$ws.Range("M98").Value = 419.0000199999999
$ws.Range("N98").Value = -4433
$ws.Range("H107").Value = 331.44446
$ws.Range("I107").Value = 298.16666
$ws.Range("J107").Value = 398
$ws.Range("K107").Value = 894.4999799999999
$ws.Range("L107").Value = 1194
$ws.Range("M107").Value = 1025.50002
$ws.Range("N107").Value = -5034
$ws.Range("H109").Value = 4950
$ws.Range("I109").Value = 4950
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 14850
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -13810
$ws.Range("H135").Value = 933.44446
$ws.Range("I135").Value = 933.44446
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8401.00014
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -5866.00014
$ws.Range("N5").ClearContents()
$ws.Range("N109").ClearContents()
$ws.Range("N135").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 18583.334
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 18583.334
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 18583.334
$ws.Range("N35").Value = -19179.334
$ws.Range("H70").Value = 6390.5454
$ws.Range("I70").Value = 3323.25
$ws.Range("J70").Value = 7072.1665
$ws.Range("K70").Value = 3323.25
$ws.Range("L70").Value = 7072.1665
$ws.Range("M70").Value = -3053.25
$ws.Range("N70").Value = -7612.1665
$ws.Range("H73").Value = 6390.5454
$ws.Range("I73").Value = 3323.25
$ws.Range("J73").Value = 7072.1665
$ws.Range("K73").Value = 3323.25
$ws.Range("L73").Value = 7072.1665
$ws.Range("M73").Value = -2387.25
$ws.Range("N73").Value = -8944.166499999999
$ws.Range("H97").Value = 2913.3333
$ws.Range("I97").Value = 2896.6
$ws.Range("K97").Value = 2896.6
$ws.Range("M97").Value = -2400.6
$ws.Range("H102").Value = 3392
$ws.Range("I102").Value = 3482.2
$ws.Range("K102").Value = 3482.2
$ws.Range("M102").Value = -1860.2
$ws.Range("H122").Value = 3626.8823
$ws.Range("I122").Value = 4023.6667
$ws.Range("J122").Value = 2674.6
$ws.Range("K122").Value = 12071.0001
$ws.Range("L122").Value = 8023.799999999999
$ws.Range("M122").Value = -9621.000100000001
$ws.Range("N122").Value = -12923.8
$ws.Range("H126").Value = 2902
$ws.Range("I126").Value = 3499.5
$ws.Range("J126").Value = 2603.25
$ws.Range("K126").Value = 10498.5
$ws.Range("L126").Value = 7809.75
$ws.Range("M126").Value = -8028.5
$ws.Range("N126").Value = -12749.75
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1523.7858
$ws.Range("I46").Value = 395
$ws.Range("J46").Value = 1711.9166
$ws.Range("K46").Value = 395
$ws.Range("L46").Value = 1711.9166
$ws.Range("M46").Value = -207
$ws.Range("N46").Value = -2087.9166
$ws.Range("H55").Value = 270.05884
$ws.Range("J55").Value = 211
$ws.Range("L55").Value = 211
$ws.Range("N55").Value = -557
$ws.Range("H68").Value = 4271.0625
$ws.Range("I68").Value = 1270.5555
$ws.Range("K68").Value = 1270.5555
$ws.Range("M68").Value = -521.5554999999999
$ws.Range("H71").Value = 4271.0625
$ws.Range("I71").Value = 1270.5555
$ws.Range("K71").Value = 6352.7775
$ws.Range("M71").Value = -2608.7775
$ws.Range("H82").Value = 3124.111
$ws.Range("J82").Value = 3693
$ws.Range("L82").Value = 3693
$ws.Range("N82").Value = -4415
$ws.Range("H85").Value = 3124.111
$ws.Range("J85").Value = 3693
$ws.Range("L85").Value = 3693
$ws.Range("N85").Value = -6189
$ws.Range("H100").Value = 2965.6667
$ws.Range("I100").Value = 2198.5
$ws.Range("K100").Value = 2198.5
$ws.Range("M100").Value = -1657.5
$ws.Range("H122").Value = 3193.6667
$ws.Range("I122").Value = 2932.6
$ws.Range("K122").Value = 8797.799999999999
$ws.Range("M122").Value = -6347.799999999999
$ws.Range("H132").Value = 3222.5
$ws.Range("I132").Value = 3222.5
$ws.Range("K132").Value = 9667.5
$ws.Range("M132").Value = -7137.5
$ws.Range("H136").Value = 50003650
$ws.Range("J136").Value = 125006250
$ws.Range("L136").Value = 375018750
$ws.Range("N136").Value = -375023850

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 55150.5
$ws.Range("J80").Value = 55150.5
$ws.Range("L80").Value = 55150.5
$ws.Range("N80").Value = -57146.5
$ws.Range("H81").Value = 5593.48
$ws.Range("I81").Value = 6509.25
$ws.Range("J81").Value = 1930.4
$ws.Range("K81").Value = 13018.5
$ws.Range("L81").Value = 3860.8
$ws.Range("M81").Value = -11957.5
$ws.Range("N81").Value = -5982.8
$ws.Range("H83").Value = 55150.5
$ws.Range("J83").Value = 55150.5
$ws.Range("L83").Value = 165451.5
$ws.Range("N83").Value = -175435.5
$ws.Range("H84").Value = 5593.48
$ws.Range("I84").Value = 6509.25
$ws.Range("J84").Value = 1930.4
$ws.Range("K84").Value = 65092.5
$ws.Range("L84").Value = 19304
$ws.Range("M84").Value = -59788.5
$ws.Range("N84").Value = -29912
$ws.Range("H99").Value = 23779.8
$ws.Range("J99").Value = 23779.8
$ws.Range("L99").Value = 23779.8
$ws.Range("N99").Value = -29769.8
$ws.Range("H122").Value = 1960.1818
$ws.Range("I122").Value = 1926.75
$ws.Range("K122").Value = 5780.25
$ws.Range("M122").Value = -3330.25
$ws.Range("H126").Value = 4998.6
$ws.Range("I126").Value = 3999
$ws.Range("J126").Value = 6498
$ws.Range("K126").Value = 11997
$ws.Range("L126").Value = 19494
$ws.Range("M126").Value = -9527
$ws.Range("N126").Value = -24434
$ws.Range("H132").Value = 2997.5
$ws.Range("I132").Value = 2997.5
$ws.Range("K132").Value = 8992.5
$ws.Range("M132").Value = -6462.5
$ws.Range("H133").Value = 24627.6
$ws.Range("J133").Value = 24627.6
$ws.Range("L133").Value = 24627.6
$ws.Range("N133").Value = -34747.6
$ws.Range("H136").Value = 1327.7142
$ws.Range("I136").Value = 1280.875
$ws.Range("J136").Value = 1477.6
$ws.Range("K136").Value = 3842.625
$ws.Range("L136").Value = 4432.799999999999
$ws.Range("M136").Value = -1292.625
$ws.Range("N136").Value = -9532.799999999999
